$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 values ---
$ws.Range("A9").Value = "Aceite de Canola"
$ws.Range("B9").Value = "Aceite"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 884
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1500

# --- A9 style: left aligned text cell, bottom-row border, not bold (like A2:A8 but last-row border) ---
$ws.Range("K1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Font.Bold = $false
$ws.Range("A9").HorizontalAlignment = -4131

# --- B9 style: left aligned, bold text cell, bottom-row border (like B2:B8 but last-row border) ---
$ws.Range("K1").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").HorizontalAlignment = -4131

# --- C9, D9, F9 style: plain percentage formatted numeric cells ---
$ws.Range("C9").NumberFormat = "0%"
$ws.Range("D9").NumberFormat = "0%"
$ws.Range("F9").NumberFormat = "0%"

# --- E9, G9, I9, K9 style: centered numeric cell, bottom-row border, not bold ---
$ws.Range("K1").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Font.Bold = $false

$ws.Range("K1").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Font.Bold = $false

$ws.Range("K1").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Font.Bold = $false

$ws.Range("K1").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Font.Bold = $false

# --- J9 style: plain integer formatted numeric cell ---
$ws.Range("J9").NumberFormat = "0"

# re-apply values lost in paste special operations
$ws.Range("A9").Value = "Aceite de Canola"
$ws.Range("B9").Value = "Aceite"
$ws.Range("E9").Value = 884
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 1500

# --- Row height ---
$ws.Rows.Item(9).RowHeight = 15.75

# --- Selection matches the saved state in the target file ---
$ws.Range("K9").Select()
